# "better gui handling of workouts"
# - Rename the workout day shown in C5 from "Tuesday" to "Monday"
# - Rename the workout group shown in C7 from "asdsads" to "Chest"
# - Add two new exercise rows (9 and 10) with exercise name / sets / reps
#   columns (C/D/E), matching the centered style already used by the
#   other label cells in this sheet
# - Widen columns C:E so the new exercise grid reads cleanly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing labels
$ws.Range("C5").Value = "Monday"
$ws.Range("C7").Value = "Chest"

# New exercise rows
$ws.Range("C9").Value = "Bench"
$ws.Range("D9").Value = "3 sets"
$ws.Range("E9").Value = "5-6 reps"

$ws.Range("C10").Value = "Incline Bench"
$ws.Range("D10").Value = "3 sets"
$ws.Range("E10").Value = "5-6 reps"

# Match the centered style used elsewhere in the sheet
$ws.Range("C9:E10").HorizontalAlignment = -4108

# Column widths for the new exercise grid (C/D/E)
$ws.Columns.Item(3).ColumnWidth = 14.1
$ws.Columns.Item(4).ColumnWidth = 6
$ws.Columns.Item(5).ColumnWidth = 8.3
